$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    [void]$d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
}

# 1. Title
Replace-Text "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING" "PRODUCT DEVELOPMENT"

# 2 & 6 & 8. "Product and Machine Learning" -> "Product Development and Product Innovation"
Replace-Text "Product and Machine Learning" "Product Development and Product Innovation"

# 3,5,7,9,12. "Product Implementation" -> "Product Development Implementation"
Replace-Text "Product Implementation" "Product Development Implementation"

# 10. ML Engineers -> Product Engineers
Replace-Text "ML Engineers" "Product Engineers"

# 11. Compliance Officers -> Quality Assurance Managers
Replace-Text "Compliance Officers" "Quality Assurance Managers"

# 13. MLflow -> Productflow
Replace-Text "MLflow" "Productflow"

# 14. Cloud ML platforms -> Cloud Product platforms
Replace-Text "Cloud ML platforms" "Cloud Product platforms"

# 15. Financial Justification: -> Product Justification:
Replace-Text "Financial Justification:" "Product Justification:"

# 16. Financial: Budget overruns... -> Product: Budget overruns...
Replace-Text "Financial: Budget overruns, cost escalation, ROI delays" "Product: Budget overruns, cost escalation, ROI delays"

# 17. ML Platform Licensing -> Product Platform Licensing
Replace-Text "ML Platform Licensing" "Product Platform Licensing"

# 18. Financial: Break-even within 30 months... -> Product: Break-even...
Replace-Text "Financial: Break-even within 30 months, 250%+ ROI within 3 years" "Product: Break-even within 30 months, 250%+ ROI within 3 years"

# 19. Financial review and budget allocation approval -> Product review and budget allocation approval
Replace-Text "Financial review and budget allocation approval" "Product review and budget allocation approval"

# 4. Remove the page break in the empty paragraph after the header block,
#    leaving an empty run behind (<w:p><w:r/></w:p>) instead of an empty
#    paragraph with no run at all.
$pageBreakPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "`f`r") {
        $pageBreakPara = $candidate
        break
    }
}
if ($pageBreakPara -ne $null) {
    $emptyRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$pageBreakPara.Range.InsertXML($emptyRunXml)
}
